$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are plain numeric-looking text (prices using "." as thousands
# separators). Force text number format first so Excel does not reinterpret them
# as numbers/dates and strip formatting (leading/trailing zeros, scientific notation).
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "30.088.04"
$ws.Range("D3").Value = "2.108.01"
$ws.Range("D5").Value = "346.10"
$ws.Range("D6").Value = "1.006"
$ws.Range("D7").Value = "0.5182"
$ws.Range("D8").Value = "0.4439"
$ws.Range("D9").Value = "0.09405"
$ws.Range("D10").Value = "52.43"
$ws.Range("D11").Value = "1.178"
$ws.Range("D12").Value = "25.41"
$ws.Range("D13").Value = "2.109.15"
$ws.Range("D14").Value = "6.754"
$ws.Range("D15").Value = "8.168"
$ws.Range("D16").Value = "100.13"
$ws.Range("D17").Value = "0.00001167"
$ws.Range("D19").Value = "20.80"
$ws.Range("D20").Value = "0.06706"
$ws.Range("D21").Value = "1.005"
$ws.Range("D22").Value = "6.241"
$ws.Range("D23").Value = "30.169.84"
$ws.Range("D24").Value = "12.72"
$ws.Range("D25").Value = "2.339"
$ws.Range("D26").Value = "2.354.31"
$ws.Range("D27").Value = "22.12"
$ws.Range("D28").Value = "2.560"
$ws.Range("D29").Value = "164.09"
$ws.Range("D30").Value = "133.94"
$ws.Range("D31").Value = "1.173"
$ws.Range("D32").Value = "0.1063"
$ws.Range("D33").Value = "1.646"
$ws.Range("D34").Value = "6.270"
$ws.Range("D35").Value = "3.955"
$ws.Range("D37").Value = "10.20"
$ws.Range("D38").Value = "0.02570"
$ws.Range("D39").Value = "0.06804"
$ws.Range("D40").Value = "0.2295"
$ws.Range("D41").Value = "0.6981"
$ws.Range("D42").Value = "12.59"
$ws.Range("D43").Value = "1.318"
$ws.Range("D44").Value = "0.6698"
$ws.Range("D45").Value = "14.27"
$ws.Range("D46").Value = "2.300"
$ws.Range("D47").Value = "3.642"
$ws.Range("D48").Value = "0.00000000354"
$ws.Range("D49").Value = "1.225"
$ws.Range("D50").Value = "82.83"
$ws.Range("D51").Value = "0.07220"

# Restore the default (unstyled) cell style so only the value itself differs,
# matching the original workbook where these cells carried no explicit style.
$dRange.Style = "Normal"

# Column E values are percentage text with surrounding whitespace padding; these
# are never numeric so no special handling is required.
$ws.Range("E2").Value = "  -2.02%  "
$ws.Range("E3").Value = "  -0.62%  "
$ws.Range("E4").Value = "  -0.92%  "
$ws.Range("E5").Value = "  +2.14%  "
$ws.Range("E6").Value = "  -0.77%  "
$ws.Range("E8").Value = "  -2.56%  "
$ws.Range("E9").Value = "  +3.05%  "
$ws.Range("E10").Value = "  -3.19%  "
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("E12").Value = "  +3.71%  "
$ws.Range("E13").Value = "  -0.66%  "
$ws.Range("E14").Value = "  -1.44%  "
$ws.Range("E15").Value = "  +0.73%  "
$ws.Range("E16").Value = "  +1.19%  "
$ws.Range("E17").Value = "  -0.46%  "
$ws.Range("E18").Value = "  -0.79%  "
$ws.Range("E19").Value = "  +6.32%  "
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("E21").Value = "  -0.79%  "
$ws.Range("E22").Value = "  -3.34%  "
$ws.Range("E23").Value = "  -1.96%  "
$ws.Range("E24").Value = "  -1.80%  "
$ws.Range("E25").Value = "  -1.63%  "
$ws.Range("E26").Value = "  -0.72%  "
$ws.Range("E27").Value = "  -1.59%  "
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("E29").Value = "  -1.09%  "
$ws.Range("E31").Value = "  -3.41%  "
$ws.Range("E32").Value = "  -1.53%  "
$ws.Range("E33").Value = "  +0.31%  "
$ws.Range("E34").Value = "  -2.39%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("E36").Value = "  +3.86%  "
$ws.Range("E37").Value = "  -3.56%  "
$ws.Range("E38").Value = "  -3.65%  "
$ws.Range("E39").Value = "  -1.25%  "
$ws.Range("E40").Value = "  -1.43%  "
$ws.Range("E41").Value = "  +0.84%  "
$ws.Range("E42").Value = "  -0.21%  "
$ws.Range("E43").Value = "  +4.00%  "
$ws.Range("E44").Value = "  +3.22%  "
$ws.Range("E45").Value = "  -6.08%  "
$ws.Range("E46").Value = "  -0.70%  "
$ws.Range("E47").Value = "  -1.76%  "
$ws.Range("E48").Value = "  -4.34%  "
$ws.Range("E49").Value = "  -2.89%  "
$ws.Range("E50").Value = "  -0.48%  "
$ws.Range("E51").Value = "  -1.33%  "
